# Update countries & provincias Spain
# Applies the data refresh captured in the commit: updated case numbers for
# several countries, swapped the "Islas Malvinas" / "Montserrat" rows back
# into alphabetical order, and bumped the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp (A1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 17 de Agosto de 2020 a las 13:58"

# --- Per-country case count updates --------------------------------------
# Row 59: Suiza
$ws.Range("B59").Value = 38252
$ws.Range("C59").Value = 128
$ws.Range("D59").Value = 33300
$ws.Range("E59").Value = 2961

# Row 70: Nepal
$ws.Range("B70").Value = 27241
$ws.Range("C70").Value = 581
$ws.Range("D70").Value = 17495
$ws.Range("E70").Value = 9639
$ws.Range("G70").Value = 3
$ws.Range("H70").Value = 107

# Row 71: Australia
$ws.Range("B71").Value = 23559
$ws.Range("C71").Value = 271
$ws.Range("D71").Value = 14539
$ws.Range("E71").Value = 8599

# Row 78: Bosnia y Herzegovina
$ws.Range("B78").Value = 16111
$ws.Range("C78").Value = 310
$ws.Range("D78").Value = 9856
$ws.Range("E78").Value = 5769
$ws.Range("G78").Value = 15
$ws.Range("H78").Value = 486

# Row 84: Sudan
$ws.Range("B84").Value = 12410
$ws.Range("C84").Value = 96
$ws.Range("D84").Value = 6385
$ws.Range("E84").Value = 5222
$ws.Range("G84").Value = 5
$ws.Range("H84").Value = 803

# Row 85: Senegal
$ws.Range("B85").Value = 12237
$ws.Range("C85").Value = 75
$ws.Range("D85").Value = 7728
$ws.Range("E85").Value = 4253
$ws.Range("G85").Value = 3
$ws.Range("H85").Value = 256

# Row 103: Croacia
$ws.Range("B103").Value = 6656
$ws.Range("C103").Value = 85
$ws.Range("D103").Value = 5254
$ws.Range("E103").Value = 1236

# Row 135: Islandia
$ws.Range("B135").Value = 2014
$ws.Range("C135").Value = 3
$ws.Range("D135").Value = 1888
$ws.Range("E135").Value = 116

# Row 158: Vietnam
$ws.Range("B158").Value = 976
$ws.Range("C158").Value = 14
$ws.Range("D158").Value = 467
$ws.Range("E158").Value = 485

# Row 174: Islas Feroe
$ws.Range("B174").Value = 373
$ws.Range("C174").Value = 1
$ws.Range("D174").Value = 229
$ws.Range("E174").Value = 144

# Row 184: Gibraltar
$ws.Range("B184").Value = 217
$ws.Range("C184").Value = 2
$ws.Range("D184").Value = 194
$ws.Range("E184").Value = 23

# --- Reorder Islas Malvinas / Montserrat back to alphabetical order ------
# Previously row 213 = Montserrat, row 214 = Islas Malvinas.
# Swap the country names and their associated data so row 213 becomes
# Islas Malvinas and row 214 becomes Montserrat (alphabetical order),
# each carrying its own stats forward.
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
